$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "331.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.26%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.50%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.721"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.97%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08301"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.41%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.799"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.90%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.500"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.97%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.972"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-5.01%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.920"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.98%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9229"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.07%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.53%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1945"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.52%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09426"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.26%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03957"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "7.36%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1067"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.27%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001313"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.08%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006038"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.54%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.512"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.98%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.72%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.119"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "10.04%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1374"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.01%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.89%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04424"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.47%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001242"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-2.50%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004443"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.81%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.67%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004000"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.17%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02815"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.06%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05616"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.00%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007953"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.62%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1426"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.60%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009054"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-4.09%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002106"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.25%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009954"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-16.03%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007372"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.33%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.17%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003610"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "11.70%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.10%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.17%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.17%"
